# Update "想去人数" (interested count) values in the F column
# of the "展览" (Exhibitions) and "全部类型" (All Types) sheets,
# matching the newly generated data from the site's scrape at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 35
$wsExhibit.Range("F4").Value = 74
$wsExhibit.Range("F5").Value = 1953
$wsExhibit.Range("F6").Value = 147
$wsExhibit.Range("F7").Value = 348

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 35
$wsAll.Range("F4").Value = 74
$wsAll.Range("F5").Value = 1953
$wsAll.Range("F6").Value = 147
$wsAll.Range("F8").Value = 348
